# Mon, Apr 27, 2020  5:04:55 AM
#
# 1) Slide 6 ("SOURCES OF FINANCE") has a table whose style is switched
#    from the custom "Table_0" style to PowerPoint's built-in
#    "Medium Style 2 - Accent 1" table style.
# 2) The deck's design ("Integral") is swapped for the default Office
#    theme colors on the slide master's theme.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{3A4C4A38-D926-4583-B1FE-8B016BC7E841}")

# --- 2. Re-colour the theme used by the slide master ----------------------
$master = $p.Designs.Item(1).SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function Set-ThemeRGB($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Office theme colours: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
Set-ThemeRGB $colorScheme 1  "000000"
Set-ThemeRGB $colorScheme 2  "FFFFFF"
Set-ThemeRGB $colorScheme 3  "44546A"
Set-ThemeRGB $colorScheme 4  "E7E6E6"
Set-ThemeRGB $colorScheme 5  "5B9BD5"
Set-ThemeRGB $colorScheme 6  "ED7D31"
Set-ThemeRGB $colorScheme 7  "A5A5A5"
Set-ThemeRGB $colorScheme 8  "FFC000"
Set-ThemeRGB $colorScheme 9  "4472C4"
Set-ThemeRGB $colorScheme 10 "70AD47"
Set-ThemeRGB $colorScheme 11 "0563C1"
Set-ThemeRGB $colorScheme 12 "954F72"
